$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.142.83'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '2.178.26'
$ws.Range('E3').Value = '  -1.96%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.54'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -7.26%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.96'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '36.34'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0935'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.18%  '
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.08%  '
$ws.Range('D15').Value = '2.506.94'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.845'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.89%  '
$ws.Range('D18').Value = '2.175.56'
$ws.Range('E18').Value = '  -2.09%  '
$ws.Range('D19').Value = '41.113.03'
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('D20').Value = '0.0₃0946'
$ws.Range('E20').Value = '  -1.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.66'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.15'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.11%  '
$ws.Range('E24').Value = '  -2.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.83'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.70%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.88%  '
$ws.Range('E28').Value = '  -4.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.85'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('E30').Value = '  -3.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.00%  '
$ws.Range('E32').Value = '  -1.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.67'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0747'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.52'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.98'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0302'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.50'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +13.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.20'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range('E42').Value = '  -7.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '60.83'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0992'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.23%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.189'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.43%  '
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('E49').Value = '  -2.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.14'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.74%  '
